# Backlog.xlsx update
# - Marks several "bancos" / report related tasks as finished ("terminado")
# - Adds two new backlog rows:
#     * "cuando se anula una OT consultar si se quieren anular la OT asociadas" (no comenzado)
#     * "facturacion en dolares " (terminado)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark finished tasks as "terminado"
$ws.Range("B9").Value  = "terminado"   # terminar circuito de movimiento de bancos
$ws.Range("B10").Value = "terminado"   # generar reporte de mov de bancos
$ws.Range("B11").Value = "terminado"   # proveedores pagos, arreglar calculo de retenciones
$ws.Range("B12").Value = "terminado"   # revisar reporte orden de pago esta fallando
$ws.Range("B15").Value = "terminado"   # reporte ot por sector ...

# New backlog items
$ws.Range("A16").Value = "cuando se anula una OT consultar si se quieren anular la OT asociadas"
$ws.Range("B16").Value = "no comenzado"

$ws.Range("A17").Value = "facturacion en dolares "
$ws.Range("B17").Value = "terminado"

# Widen column A slightly to keep the "best fit" look after the new, longer rows
$ws.Columns.Item(1).ColumnWidth = 73.83

# Leave the selection where the user would naturally continue typing
$ws.Range("B18").Select()

$wb.Save()
